$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''68.358.98'
$ws.Range('E2').Value = '''  -1.64%  '
$ws.Range('D3').Value = '''2.435.34'
$ws.Range('E3').Value = '''  -2.19%  '
$ws.Range('E4').Value = '''  -0.02%  '
$ws.Range('D5').Value = '''553.82'
$ws.Range('E5').Value = '''  -2.71%  '
$ws.Range('D6').Value = '''160.01'
$ws.Range('E6').Value = '''  -2.64%  '
$ws.Range('E7').Value = '''  +0.03%  '
$ws.Range('E8').Value = '''  -2.35%  '
$ws.Range('D9').Value = '''2.433.30'
$ws.Range('E9').Value = '''  -2.20%  '
$ws.Range('D10').Value = '''0.147'
$ws.Range('E10').Value = '''  -7.98%  '
$ws.Range('E11').Value = '''  -1.09%  '
$ws.Range('D12').Value = '''0.331'
$ws.Range('E12').Value = '''  -5.76%  '
$ws.Range('E13').Value = '''  -2.38%  '
$ws.Range('D14').Value = '''2.880.15'
$ws.Range('E14').Value = '''  -2.18%  '
$ws.Range('D15').Value = '''68.226.07'
$ws.Range('E15').Value = '''  -1.58%  '
$ws.Range('D16').Value = '''0.0000166'
$ws.Range('E16').Value = '''  -5.01%  '
$ws.Range('D17').Value = '''23.13'
$ws.Range('E17').Value = '''  -4.11%  '
$ws.Range('D18').Value = '''2.431.85'
$ws.Range('E18').Value = '''  -2.09%  '
$ws.Range('D19').Value = '''10.64'
$ws.Range('E19').Value = '''  -4.72%  '
$ws.Range('D20').Value = '''337.35'
$ws.Range('E20').Value = '''  -2.72%  '
$ws.Range('D22').Value = '''3.76'
$ws.Range('E22').Value = '''  -2.94%  '
$ws.Range('B23').Value = '''LEO'
$ws.Range('C23').Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D23').Value = '''5.96'
$ws.Range('E23').Value = '''  -1.84%  '
$ws.Range('B24').Value = '''Dai'
$ws.Range('C24').Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '''1.00'
$ws.Range('E24').Value = '''  +0.03%  '
$ws.Range('B25').Value = '''SuiNetwork'
$ws.Range('C25').Value = '''https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').Value = '''1.84'
$ws.Range('E25').Value = '''  -2.74%  '
$ws.Range('B26').Value = '''Litecoin'
$ws.Range('C26').Value = '''https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '''65.89'
$ws.Range('E26').Value = '''  -4.82%  '
$ws.Range('B27').Value = '''WrappedeETH'
$ws.Range('C27').Value = '''https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '''2.566.33'
$ws.Range('E27').Value = '''  -1.91%  '
$ws.Range('B28').Value = '''NEARProtocol'
$ws.Range('C28').Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D28').Value = '''3.60'
$ws.Range('E28').Value = '''  -6.86%  '
$ws.Range('B29').Value = '''Binance-PegBSC-USD'
$ws.Range('C29').Value = '''https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '''0.997'
$ws.Range('E29').Value = '''  -0.29%  '
$ws.Range('B30').Value = '''Aptos'
$ws.Range('C30').Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '''7.99'
$ws.Range('E30').Value = '''  -6.80%  '
$ws.Range('B31').Value = '''PEPE'
$ws.Range('C31').Value = '''https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '''0.0₃0802'
$ws.Range('E31').Value = '''  -7.54%  '
$ws.Range('B32').Value = '''InternetComputer(DFINITY)'
$ws.Range('C32').Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''7.07'
$ws.Range('E32').Value = '''  -6.39%  '
$ws.Range('B33').Value = '''FirstDigitalUSD'
$ws.Range('C33').Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').Value = '''0.999'
$ws.Range('E33').Value = '''  -0.02%  '
$ws.Range('B34').Value = '''Bittensor'
$ws.Range('C34').Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = '''429.34'
$ws.Range('E34').Value = '''  -1.20%  '
$ws.Range('B35').Value = '''Fetch.AI'
$ws.Range('C35').Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '''1.11'
$ws.Range('E35').Value = '''  -6.42%  '
$ws.Range('B36').Value = '''PancakeSwap'
$ws.Range('C36').Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D36').Value = '''1.59'
$ws.Range('E36').Value = '''  -6.29%  '
$ws.Range('B37').Value = '''Monero'
$ws.Range('C37').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '''155.59'
$ws.Range('E37').Value = '''  +0.17%  '
$ws.Range('B38').Value = '''WhiteBITCoin'
$ws.Range('C38').Value = '''https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').Value = '''19.02'
$ws.Range('E38').Value = '''  -0.28%  '
$ws.Range('B39').Value = '''USDe'
$ws.Range('C39').Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '''  -0.05%  '
$ws.Range('B40').Value = '''Kaspa'
$ws.Range('C40').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '''0.108'
$ws.Range('E40').Value = '''  -3.57%  '
$ws.Range('B41').Value = '''EthereumClassic'
$ws.Range('C41').Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').Value = '''17.66'
$ws.Range('E41').Value = '''  -2.39%  '
$ws.Range('B42').Value = '''PolygonEcosystemToken'
$ws.Range('C42').Value = '''https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').Value = '''0.300'
$ws.Range('E42').Value = '''  -4.14%  '
$ws.Range('B43').Value = '''RenderToken'
$ws.Range('C43').Value = '''https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').Value = '''4.36'
$ws.Range('E43').Value = '''  -4.73%  '
$ws.Range('B44').Value = '''OKB'
$ws.Range('C44').Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '''37.35'
$ws.Range('E44').Value = '''  -1.01%  '
$ws.Range('B45').Value = '''Stacks'
$ws.Range('C45').Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '''1.44'
$ws.Range('E45').Value = '''  -8.50%  '
$ws.Range('B46').Value = '''ImmutableX'
$ws.Range('C46').Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').Value = '''1.07'
$ws.Range('E46').Value = '''  +0.35%  '
$ws.Range('B47').Value = '''dogwifhat'
$ws.Range('C47').Value = '''https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '''2.00'
$ws.Range('E47').Value = '''  -7.18%  '
$ws.Range('B48').Value = '''Aave'
$ws.Range('C48').Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '''130.79'
$ws.Range('E48').Value = '''  -5.49%  '
$ws.Range('B49').Value = '''Filecoin'
$ws.Range('C49').Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = '''3.30'
$ws.Range('E49').Value = '''  -3.35%  '
$ws.Range('B50').Value = '''Cronos'
$ws.Range('C50').Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.0711'
$ws.Range('E50').Value = '''  -1.49%  '
$ws.Range('B51').Value = '''ARBITRUM'
$ws.Range('C51').Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = '''0.476'
$ws.Range('E51').Value = '''  -5.54%  '
